# Excel COM-interop script applying commit "Update latest output (run 221)"
# to the optimisation_result workbook: Schedule gets 2 more scheduled pump
# runs (and a recompute of row 3's cost columns), Detailed gets 35 rows
# reclassified from forecast->historical (with revised actuals) plus a full
# new day (rows 50-97) of forecast data for 2026-02-07.

$wb = $excel.ActiveWorkbook

# ----- Sheet: Schedule -----
$schedule = $wb.Worksheets.Item("Schedule")

# Row 3's Cost/Unit Cost were recomputed with the new optimisation run
$schedule.Cells.Item(3, 5).Value = 275.8933807500001
$schedule.Cells.Item(3, 6).Value = 10.4268095521542

# Append newly scheduled pump runs (rows 4 and 5)
$schedule.Cells.Item(4, 1).Value = 46060.27083333334
$schedule.Cells.Item(4, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(4, 2).Value = 46060.5625
$schedule.Cells.Item(4, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(4, 3).Value = 7
$schedule.Cells.Item(4, 4).Value = 26.46
$schedule.Cells.Item(4, 5).Value = 935.2722405
$schedule.Cells.Item(4, 6).Value = 35.34664552154195

$schedule.Cells.Item(5, 1).Value = 46060.83333333334
$schedule.Cells.Item(5, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(5, 2).Value = 46061
$schedule.Cells.Item(5, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(5, 3).Value = 4
$schedule.Cells.Item(5, 4).Value = 15.12
$schedule.Cells.Item(5, 5).Value = 746.17798125
$schedule.Cells.Item(5, 6).Value = 49.35039558531746

# ----- Sheet: Detailed -----
$detailed = $wb.Worksheets.Item("Detailed")

# Re-run of the forecast/historical split for 2026-02-06 (rows 14-48):
# later actuals replaced several forecast rows with historical readings.
$detailed.Cells.Item(14, 2).Value = 286.19889
$detailed.Cells.Item(15, 2).Value = 245.61799
$detailed.Cells.Item(16, 3).Value = "historical"
$detailed.Cells.Item(17, 2).Value = 61.03948
$detailed.Cells.Item(17, 3).Value = "historical"
$detailed.Cells.Item(18, 2).Value = 59.15131
$detailed.Cells.Item(18, 3).Value = "historical"
$detailed.Cells.Item(19, 2).Value = 56.98
$detailed.Cells.Item(19, 3).Value = "historical"
$detailed.Cells.Item(20, 3).Value = "historical"
$detailed.Cells.Item(21, 2).Value = 36.06
$detailed.Cells.Item(21, 3).Value = "historical"
$detailed.Cells.Item(22, 3).Value = "historical"
$detailed.Cells.Item(23, 2).Value = 0.50992
$detailed.Cells.Item(23, 3).Value = "historical"
$detailed.Cells.Item(24, 2).Value = 0.50993
$detailed.Cells.Item(24, 3).Value = "historical"
$detailed.Cells.Item(25, 2).Value = 0.50993
$detailed.Cells.Item(25, 3).Value = "historical"
$detailed.Cells.Item(26, 2).Value = -0.13867
$detailed.Cells.Item(26, 3).Value = "historical"
$detailed.Cells.Item(27, 2).Value = -5.50985
$detailed.Cells.Item(27, 3).Value = "historical"
$detailed.Cells.Item(28, 3).Value = "historical"
$detailed.Cells.Item(29, 2).Value = 36.06
$detailed.Cells.Item(29, 3).Value = "historical"
$detailed.Cells.Item(30, 2).Value = 0.71552
$detailed.Cells.Item(30, 3).Value = "historical"
$detailed.Cells.Item(31, 3).Value = "historical"
$detailed.Cells.Item(32, 2).Value = 153.01
$detailed.Cells.Item(32, 3).Value = "historical"
$detailed.Cells.Item(33, 2).Value = 153.01
$detailed.Cells.Item(33, 3).Value = "historical"
$detailed.Cells.Item(34, 2).Value = 5300.46487
$detailed.Cells.Item(35, 2).Value = 600.0
$detailed.Cells.Item(36, 2).Value = 153.01
$detailed.Cells.Item(37, 2).Value = 189.13
$detailed.Cells.Item(38, 2).Value = 162.52339
$detailed.Cells.Item(39, 2).Value = 216.29664
$detailed.Cells.Item(40, 2).Value = 152.36792
$detailed.Cells.Item(41, 2).Value = 174.03554
$detailed.Cells.Item(42, 2).Value = 299.75
$detailed.Cells.Item(43, 2).Value = 189.85
$detailed.Cells.Item(44, 2).Value = 182.05854
$detailed.Cells.Item(45, 2).Value = 149.86851
$detailed.Cells.Item(46, 2).Value = 146.14775
$detailed.Cells.Item(47, 2).Value = 158.91093
$detailed.Cells.Item(48, 2).Value = 138.42

# Append newly forecast period for 2026-02-07 (rows 50-97)
$detailed.Cells.Item(50, 1).Value = 46060.0
$detailed.Cells.Item(50, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(50, 2).Value = 166.09405
$detailed.Cells.Item(50, 3).Value = "forecast"
$detailed.Cells.Item(50, 4).Value = 46060.0
$detailed.Cells.Item(50, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(50, 5).Value = "OFF"

$detailed.Cells.Item(51, 1).Value = 46060.02083333334
$detailed.Cells.Item(51, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(51, 2).Value = 138.42
$detailed.Cells.Item(51, 3).Value = "forecast"
$detailed.Cells.Item(51, 4).Value = 46060.0
$detailed.Cells.Item(51, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(51, 5).Value = "OFF"

$detailed.Cells.Item(52, 1).Value = 46060.04166666666
$detailed.Cells.Item(52, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(52, 2).Value = 138.42
$detailed.Cells.Item(52, 3).Value = "forecast"
$detailed.Cells.Item(52, 4).Value = 46060.0
$detailed.Cells.Item(52, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(52, 5).Value = "OFF"

$detailed.Cells.Item(53, 1).Value = 46060.0625
$detailed.Cells.Item(53, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(53, 2).Value = 114.58672
$detailed.Cells.Item(53, 3).Value = "forecast"
$detailed.Cells.Item(53, 4).Value = 46060.0
$detailed.Cells.Item(53, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(53, 5).Value = "OFF"

$detailed.Cells.Item(54, 1).Value = 46060.08333333334
$detailed.Cells.Item(54, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(54, 2).Value = 119.84965
$detailed.Cells.Item(54, 3).Value = "forecast"
$detailed.Cells.Item(54, 4).Value = 46060.0
$detailed.Cells.Item(54, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(54, 5).Value = "OFF"

$detailed.Cells.Item(55, 1).Value = 46060.10416666666
$detailed.Cells.Item(55, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(55, 2).Value = 114.74563
$detailed.Cells.Item(55, 3).Value = "forecast"
$detailed.Cells.Item(55, 4).Value = 46060.0
$detailed.Cells.Item(55, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(55, 5).Value = "OFF"

$detailed.Cells.Item(56, 1).Value = 46060.125
$detailed.Cells.Item(56, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(56, 2).Value = 116.46553
$detailed.Cells.Item(56, 3).Value = "forecast"
$detailed.Cells.Item(56, 4).Value = 46060.0
$detailed.Cells.Item(56, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(56, 5).Value = "OFF"

$detailed.Cells.Item(57, 1).Value = 46060.14583333334
$detailed.Cells.Item(57, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(57, 2).Value = 115.67607
$detailed.Cells.Item(57, 3).Value = "forecast"
$detailed.Cells.Item(57, 4).Value = 46060.0
$detailed.Cells.Item(57, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(57, 5).Value = "OFF"

$detailed.Cells.Item(58, 1).Value = 46060.16666666666
$detailed.Cells.Item(58, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(58, 2).Value = 115.7191
$detailed.Cells.Item(58, 3).Value = "forecast"
$detailed.Cells.Item(58, 4).Value = 46060.0
$detailed.Cells.Item(58, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(58, 5).Value = "OFF"

$detailed.Cells.Item(59, 1).Value = 46060.1875
$detailed.Cells.Item(59, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(59, 2).Value = 108.89
$detailed.Cells.Item(59, 3).Value = "forecast"
$detailed.Cells.Item(59, 4).Value = 46060.0
$detailed.Cells.Item(59, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(59, 5).Value = "OFF"

$detailed.Cells.Item(60, 1).Value = 46060.20833333334
$detailed.Cells.Item(60, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(60, 2).Value = 114.16346
$detailed.Cells.Item(60, 3).Value = "forecast"
$detailed.Cells.Item(60, 4).Value = 46060.0
$detailed.Cells.Item(60, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(60, 5).Value = "OFF"

$detailed.Cells.Item(61, 1).Value = 46060.22916666666
$detailed.Cells.Item(61, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(61, 2).Value = 127.68089
$detailed.Cells.Item(61, 3).Value = "forecast"
$detailed.Cells.Item(61, 4).Value = 46060.0
$detailed.Cells.Item(61, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(61, 5).Value = "OFF"

$detailed.Cells.Item(62, 1).Value = 46060.25
$detailed.Cells.Item(62, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(62, 2).Value = 120.9633
$detailed.Cells.Item(62, 3).Value = "forecast"
$detailed.Cells.Item(62, 4).Value = 46060.0
$detailed.Cells.Item(62, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(62, 5).Value = "OFF"

$detailed.Cells.Item(63, 1).Value = 46060.27083333334
$detailed.Cells.Item(63, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(63, 2).Value = 89.65441
$detailed.Cells.Item(63, 3).Value = "forecast"
$detailed.Cells.Item(63, 4).Value = 46060.0
$detailed.Cells.Item(63, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(63, 5).Value = "ON"

$detailed.Cells.Item(64, 1).Value = 46060.29166666666
$detailed.Cells.Item(64, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(64, 2).Value = 65.00007
$detailed.Cells.Item(64, 3).Value = "forecast"
$detailed.Cells.Item(64, 4).Value = 46060.0
$detailed.Cells.Item(64, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(64, 5).Value = "ON"

$detailed.Cells.Item(65, 1).Value = 46060.3125
$detailed.Cells.Item(65, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(65, 2).Value = 63.56936
$detailed.Cells.Item(65, 3).Value = "forecast"
$detailed.Cells.Item(65, 4).Value = 46060.0
$detailed.Cells.Item(65, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(65, 5).Value = "ON"

$detailed.Cells.Item(66, 1).Value = 46060.33333333334
$detailed.Cells.Item(66, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(66, 2).Value = 57.08
$detailed.Cells.Item(66, 3).Value = "forecast"
$detailed.Cells.Item(66, 4).Value = 46060.0
$detailed.Cells.Item(66, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(66, 5).Value = "ON"

$detailed.Cells.Item(67, 1).Value = 46060.35416666666
$detailed.Cells.Item(67, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(67, 2).Value = 57.08
$detailed.Cells.Item(67, 3).Value = "forecast"
$detailed.Cells.Item(67, 4).Value = 46060.0
$detailed.Cells.Item(67, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(67, 5).Value = "ON"

$detailed.Cells.Item(68, 1).Value = 46060.375
$detailed.Cells.Item(68, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(68, 2).Value = 36.0595
$detailed.Cells.Item(68, 3).Value = "forecast"
$detailed.Cells.Item(68, 4).Value = 46060.0
$detailed.Cells.Item(68, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(68, 5).Value = "ON"

$detailed.Cells.Item(69, 1).Value = 46060.39583333334
$detailed.Cells.Item(69, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(69, 2).Value = 56.98
$detailed.Cells.Item(69, 3).Value = "forecast"
$detailed.Cells.Item(69, 4).Value = 46060.0
$detailed.Cells.Item(69, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(69, 5).Value = "ON"

$detailed.Cells.Item(70, 1).Value = 46060.41666666666
$detailed.Cells.Item(70, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(70, 2).Value = 56.98
$detailed.Cells.Item(70, 3).Value = "forecast"
$detailed.Cells.Item(70, 4).Value = 46060.0
$detailed.Cells.Item(70, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(70, 5).Value = "ON"

$detailed.Cells.Item(71, 1).Value = 46060.4375
$detailed.Cells.Item(71, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(71, 2).Value = 36.0601
$detailed.Cells.Item(71, 3).Value = "forecast"
$detailed.Cells.Item(71, 4).Value = 46060.0
$detailed.Cells.Item(71, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(71, 5).Value = "ON"

$detailed.Cells.Item(72, 1).Value = 46060.45833333334
$detailed.Cells.Item(72, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(72, 2).Value = 57.06007
$detailed.Cells.Item(72, 3).Value = "forecast"
$detailed.Cells.Item(72, 4).Value = 46060.0
$detailed.Cells.Item(72, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(72, 5).Value = "ON"

$detailed.Cells.Item(73, 1).Value = 46060.47916666666
$detailed.Cells.Item(73, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(73, 2).Value = 57.06007
$detailed.Cells.Item(73, 3).Value = "forecast"
$detailed.Cells.Item(73, 4).Value = 46060.0
$detailed.Cells.Item(73, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(73, 5).Value = "ON"

$detailed.Cells.Item(74, 1).Value = 46060.5
$detailed.Cells.Item(74, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(74, 2).Value = 108.89
$detailed.Cells.Item(74, 3).Value = "forecast"
$detailed.Cells.Item(74, 4).Value = 46060.0
$detailed.Cells.Item(74, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(74, 5).Value = "ON"

$detailed.Cells.Item(75, 1).Value = 46060.52083333334
$detailed.Cells.Item(75, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(75, 2).Value = 108.89
$detailed.Cells.Item(75, 3).Value = "forecast"
$detailed.Cells.Item(75, 4).Value = 46060.0
$detailed.Cells.Item(75, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(75, 5).Value = "ON"

$detailed.Cells.Item(76, 1).Value = 46060.54166666666
$detailed.Cells.Item(76, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(76, 2).Value = 108.89
$detailed.Cells.Item(76, 3).Value = "forecast"
$detailed.Cells.Item(76, 4).Value = 46060.0
$detailed.Cells.Item(76, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(76, 5).Value = "ON"

$detailed.Cells.Item(77, 1).Value = 46060.5625
$detailed.Cells.Item(77, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(77, 2).Value = 599.99
$detailed.Cells.Item(77, 3).Value = "forecast"
$detailed.Cells.Item(77, 4).Value = 46060.0
$detailed.Cells.Item(77, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(77, 5).Value = "OFF"

$detailed.Cells.Item(78, 1).Value = 46060.58333333334
$detailed.Cells.Item(78, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(78, 2).Value = 12532.66789
$detailed.Cells.Item(78, 3).Value = "forecast"
$detailed.Cells.Item(78, 4).Value = 46060.0
$detailed.Cells.Item(78, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(78, 5).Value = "OFF"

$detailed.Cells.Item(79, 1).Value = 46060.60416666666
$detailed.Cells.Item(79, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(79, 2).Value = 20203.88
$detailed.Cells.Item(79, 3).Value = "forecast"
$detailed.Cells.Item(79, 4).Value = 46060.0
$detailed.Cells.Item(79, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(79, 5).Value = "OFF"

$detailed.Cells.Item(80, 1).Value = 46060.625
$detailed.Cells.Item(80, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(80, 2).Value = 20234.5
$detailed.Cells.Item(80, 3).Value = "forecast"
$detailed.Cells.Item(80, 4).Value = 46060.0
$detailed.Cells.Item(80, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(80, 5).Value = "OFF"

$detailed.Cells.Item(81, 1).Value = 46060.64583333334
$detailed.Cells.Item(81, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(81, 2).Value = 19036.44
$detailed.Cells.Item(81, 3).Value = "forecast"
$detailed.Cells.Item(81, 4).Value = 46060.0
$detailed.Cells.Item(81, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(81, 5).Value = "OFF"

$detailed.Cells.Item(82, 1).Value = 46060.66666666666
$detailed.Cells.Item(82, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(82, 2).Value = 15019.16219
$detailed.Cells.Item(82, 3).Value = "forecast"
$detailed.Cells.Item(82, 4).Value = 46060.0
$detailed.Cells.Item(82, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(82, 5).Value = "OFF"

$detailed.Cells.Item(83, 1).Value = 46060.6875
$detailed.Cells.Item(83, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(83, 2).Value = 20234.5
$detailed.Cells.Item(83, 3).Value = "forecast"
$detailed.Cells.Item(83, 4).Value = 46060.0
$detailed.Cells.Item(83, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(83, 5).Value = "OFF"

$detailed.Cells.Item(84, 1).Value = 46060.70833333334
$detailed.Cells.Item(84, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(84, 2).Value = 14943.11133
$detailed.Cells.Item(84, 3).Value = "forecast"
$detailed.Cells.Item(84, 4).Value = 46060.0
$detailed.Cells.Item(84, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(84, 5).Value = "OFF"

$detailed.Cells.Item(85, 1).Value = 46060.72916666666
$detailed.Cells.Item(85, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(85, 2).Value = 12387.732
$detailed.Cells.Item(85, 3).Value = "forecast"
$detailed.Cells.Item(85, 4).Value = 46060.0
$detailed.Cells.Item(85, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(85, 5).Value = "OFF"

$detailed.Cells.Item(86, 1).Value = 46060.75
$detailed.Cells.Item(86, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(86, 2).Value = 12390.1261
$detailed.Cells.Item(86, 3).Value = "forecast"
$detailed.Cells.Item(86, 4).Value = 46060.0
$detailed.Cells.Item(86, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(86, 5).Value = "OFF"

$detailed.Cells.Item(87, 1).Value = 46060.77083333334
$detailed.Cells.Item(87, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(87, 2).Value = 12412.75052
$detailed.Cells.Item(87, 3).Value = "forecast"
$detailed.Cells.Item(87, 4).Value = 46060.0
$detailed.Cells.Item(87, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(87, 5).Value = "OFF"

$detailed.Cells.Item(88, 1).Value = 46060.79166666666
$detailed.Cells.Item(88, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(88, 2).Value = 274.32458
$detailed.Cells.Item(88, 3).Value = "forecast"
$detailed.Cells.Item(88, 4).Value = 46060.0
$detailed.Cells.Item(88, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(88, 5).Value = "OFF"

$detailed.Cells.Item(89, 1).Value = 46060.8125
$detailed.Cells.Item(89, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(89, 2).Value = 156.34441
$detailed.Cells.Item(89, 3).Value = "forecast"
$detailed.Cells.Item(89, 4).Value = 46060.0
$detailed.Cells.Item(89, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(89, 5).Value = "OFF"

$detailed.Cells.Item(90, 1).Value = 46060.83333333334
$detailed.Cells.Item(90, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(90, 2).Value = 95.81988
$detailed.Cells.Item(90, 3).Value = "forecast"
$detailed.Cells.Item(90, 4).Value = 46060.0
$detailed.Cells.Item(90, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(90, 5).Value = "ON"

$detailed.Cells.Item(91, 1).Value = 46060.85416666666
$detailed.Cells.Item(91, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(91, 2).Value = 100.8224
$detailed.Cells.Item(91, 3).Value = "forecast"
$detailed.Cells.Item(91, 4).Value = 46060.0
$detailed.Cells.Item(91, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(91, 5).Value = "ON"

$detailed.Cells.Item(92, 1).Value = 46060.875
$detailed.Cells.Item(92, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(92, 2).Value = 75.89678
$detailed.Cells.Item(92, 3).Value = "forecast"
$detailed.Cells.Item(92, 4).Value = 46060.0
$detailed.Cells.Item(92, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(92, 5).Value = "ON"

$detailed.Cells.Item(93, 1).Value = 46060.89583333334
$detailed.Cells.Item(93, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(93, 2).Value = 108.89
$detailed.Cells.Item(93, 3).Value = "forecast"
$detailed.Cells.Item(93, 4).Value = 46060.0
$detailed.Cells.Item(93, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(93, 5).Value = "ON"

$detailed.Cells.Item(94, 1).Value = 46060.91666666666
$detailed.Cells.Item(94, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(94, 2).Value = 108.89
$detailed.Cells.Item(94, 3).Value = "forecast"
$detailed.Cells.Item(94, 4).Value = 46060.0
$detailed.Cells.Item(94, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(94, 5).Value = "ON"

$detailed.Cells.Item(95, 1).Value = 46060.9375
$detailed.Cells.Item(95, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(95, 2).Value = 89.07312
$detailed.Cells.Item(95, 3).Value = "forecast"
$detailed.Cells.Item(95, 4).Value = 46060.0
$detailed.Cells.Item(95, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(95, 5).Value = "ON"

$detailed.Cells.Item(96, 1).Value = 46060.95833333334
$detailed.Cells.Item(96, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(96, 2).Value = 93.09766
$detailed.Cells.Item(96, 3).Value = "forecast"
$detailed.Cells.Item(96, 4).Value = 46060.0
$detailed.Cells.Item(96, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(96, 5).Value = "ON"

$detailed.Cells.Item(97, 1).Value = 46060.97916666666
$detailed.Cells.Item(97, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(97, 2).Value = 92.82091
$detailed.Cells.Item(97, 3).Value = "forecast"
$detailed.Cells.Item(97, 4).Value = 46060.0
$detailed.Cells.Item(97, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(97, 5).Value = "ON"

Write-Output "Applied run 221 update: Schedule +2 rows; Detailed +48 rows, 35 reclassified, 2 recomputed totals"
